$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, shifting existing rows 12-40 down to 13-41
$ws.Rows.Item(12).Insert()

# Fill the new row 12 with the "Introduction to SQL" course entry
$ws.Cells.Item(12, 1).Value = "Introduction to SQL"
$ws.Cells.Item(12, 7).Value = 3

# Update the selection / view to match the post-edit state
$ws.Range("G13").Select()
$ws.Application.ActiveWindow.ScrollRow = 4
